$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reposition the workbook window (best effort; matches the window move in the diff)
$win = $wb.Windows.Item(1)
$win.Left = 11180
$win.Top = 3240

# Add new row 19 data: seriesId mutability info
$ws.Range("A19").Value = "seriesId"
$ws.Range("B19").Value = "Client"
$ws.Range("C19").Value = "Client"
$ws.Range("D19").Value = "No"

# Update selection to match the new active cell shown in the diff
$ws.Range("B20").Select() | Out-Null
